$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = 2000
$ws.Range("B22").Value = "struggle"
$ws.Range("C22").Value = -5.779280513524995
$ws.Range("D22").Value = -3.271673738956447
$ws.Range("E22").Value = 4.346601516008363
$ws.Range("F22").Value = -0.2765692472457886
$ws.Range("G22").Value = 2.626567840576172
$ws.Range("H22").Value = 0.6087272167205811

$ws.Range("A23").Value = 2100
$ws.Range("B23").Value = "struggle"
$ws.Range("C23").Value = 3.140387788414934
$ws.Range("D23").Value = 0.259726375341407
$ws.Range("E23").Value = -4.785581156611421
$ws.Range("F23").Value = -0.3174972236156463
$ws.Range("G23").Value = 1.27487576007843
$ws.Range("H23").Value = 0.1149953827261924

$ws.Range("A24").Value = 2200
$ws.Range("B24").Value = "struggle"
$ws.Range("C24").Value = -1.632258296012878
$ws.Range("D24").Value = 0.6425724923610687
$ws.Range("E24").Value = -3.22618693113327
$ws.Range("F24").Value = -0.2593123018741607
$ws.Range("G24").Value = -0.741895854473114
$ws.Range("H24").Value = -0.290313720703125

$ws.Range("A25").Value = 2300
$ws.Range("B25").Value = "struggle"
$ws.Range("C25").Value = -3.271841421723368
$ws.Range("D25").Value = 0.07577018067240526
$ws.Range("E25").Value = -1.077775649726385
$ws.Range("F25").Value = 0.6250678896903992
$ws.Range("G25").Value = -3.214983701705933
$ws.Range("H25").Value = -0.8246681094169617

$ws.Range("A26").Value = 2400
$ws.Range("B26").Value = "struggle"
$ws.Range("C26").Value = -3.740465611219407
$ws.Range("D26").Value = -0.2502757757902145
$ws.Range("E26").Value = -5.230584308505059
$ws.Range("F26").Value = 0.2973386645317077
$ws.Range("G26").Value = -4.603633403778076
$ws.Range("H26").Value = 0.0426078513264656

$ws.Range("A27").Value = 2500
$ws.Range("B27").Value = "struggle"
$ws.Range("C27").Value = -1.173786669969556
$ws.Range("D27").Value = -1.206141140311958
$ws.Range("E27").Value = -5.999948702752588
$ws.Range("F27").Value = -0.3110831379890442
$ws.Range("G27").Value = -4.220009803771973
$ws.Range("H27").Value = 1.416138410568237

$ws.Range("A28").Value = 2600
$ws.Range("B28").Value = "struggle"
$ws.Range("C28").Value = -0.527452439069747
$ws.Range("D28").Value = -1.933494433760643
$ws.Range("E28").Value = -1.995455801486972
$ws.Range("F28").Value = -0.2417499274015426
$ws.Range("G28").Value = -3.320205211639404
$ws.Range("H28").Value = 1.446528911590576

$ws.Range("A29").Value = 2700
$ws.Range("B29").Value = "struggle"
$ws.Range("C29").Value = -2.818732134997842
$ws.Range("D29").Value = -1.367858927696945
$ws.Range("E29").Value = 0.861171409487724
$ws.Range("F29").Value = 0.3194825351238251
$ws.Range("G29").Value = -0.44180828332901
$ws.Range("H29").Value = 0.5250386595726013

$ws.Range("A30").Value = 2800
$ws.Range("B30").Value = "struggle"
$ws.Range("C30").Value = -4.338251754641532
$ws.Range("D30").Value = -0.358771674335002
$ws.Range("E30").Value = 0.811524987220763
$ws.Range("F30").Value = 0.0740674138069152
$ws.Range("G30").Value = 2.8290696144104
$ws.Range("H30").Value = -0.7684684991836548

$ws.Range("A31").Value = 2900
$ws.Range("B31").Value = "struggle"
$ws.Range("C31").Value = -5.003720842301852
$ws.Range("D31").Value = -0.2583636995404971
$ws.Range("E31").Value = 1.431181490421301
$ws.Range("F31").Value = 0.3081815242767334
$ws.Range("G31").Value = 5.233893394470215
$ws.Range("H31").Value = -0.9094256162643432

